# Generate Report for Handback
# Update the timestamp strings recorded in the handback status workbook.

$wb = $excel.ActiveWorkbook

# Overview sheet: "Latest HO Xliff Generate Date" for the
# ba303ce5-72b3-4ebb-a135-ceb5b863046b.md row.
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G4").Value = "2016-08-25 04:44:02"

# zh-cn sheet: handoff / handback datetimes for the same file's row.
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("H4").Value = "2016-08-25 04:43:56"
$wsZhCn.Range("K4").Value = "2016-08-25 04:44:28"

# de-de sheet: matching "Latest HO Xliff Generate Date" plus the
# handback datetime for the same row.
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("H4").Value = "2016-08-25 04:44:02"
$wsDeDe.Range("K4").Value = "2016-08-25 04:44:35"
